$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new headers for age and ageband in columns H and I
$ws.Range("H1").Value = "age"
$ws.Range("I1").Value = "ageband"

# Rename column header C1 from "sex" to "gender"
$ws.Range("C1").Value = "gender"

# Keep final selection on C1, matching the recorded edit state
$ws.Range("C1").Select()
